# Update sheet title and data for the new "through 03-19" figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name shown in workbook.xml <sheets><sheet name=.../>)
$ws.Name = "Through 2022-03-19"

# Update the row label for March
$ws.Range("A4").Value = "March (through 03-19)"

# Row 4: March data, columns C:I (2016-2022)
$ws.Range("C4").Value = 28
$ws.Range("D4").Value = 34
$ws.Range("E4").Value = 38
$ws.Range("F4").Value = 20
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 84

# Row 5: Total data, columns C:I (2016-2022)
$ws.Range("C5").Value = 115
$ws.Range("D5").Value = 165
$ws.Range("E5").Value = 175
$ws.Range("F5").Value = 99
$ws.Range("G5").Value = 181
$ws.Range("H5").Value = 392
$ws.Range("I5").Value = 384
